# Work Profile and new tenant support
# Appends new sprint-run history rows to the AMSIN, BETA and AMS sheets,
# and fixes up the formatting/precision of AMS!B41 (which previously had
# no explicit style and a slightly different cached float).
#
# NOTE: this host's PowerShell function calls only bind POSITIONAL
# arguments correctly (named "-Param Value" args come through empty), so
# every helper below is called positionally.

$wb = $excel.ActiveWorkbook

function Add-HistoryRow($ws, $row, $dateText, $runTime, $sprintName, $total, $pass, $fail, $timeTaken, $formatSourceRow) {
    # Column A (Run Date) - stored as literal text like "2023-03-10", not
    # an actual date serial, so force text entry then drop back to a
    # plain/Normal style (keeps the same look as the rest of the column,
    # which uses the sheet's default "General" style).
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $dateText
    $ws.Cells.Item($row, 1).Style = "Normal"

    # Column B (Run Time) - numeric date/time serial, formatted the same
    # way as the row above it.
    $ws.Cells.Item($formatSourceRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $runTime

    # Column C (Sprint Name) - also literal text, same trick as column A.
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 3).Style = "Normal"

    # Columns D-G (Total/Pass/Fail Cases, Time Taken) - plain numbers using
    # the column/sheet default style.
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN sheet: add rows 76-79
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-HistoryRow $wsAmsin 76 "2023-03-10" 44995.8066925926 "174ffiinnalrun" 89 87 2 2.88 75
Add-HistoryRow $wsAmsin 77 "2023-03-13" 44998.48121105324 "174finalrun" 89 87 2 2.28 76
Add-HistoryRow $wsAmsin 78 "2023-03-31" 45016.45243590278 "175fnlrun" 89 87 2 2.55 77

# Row 79 is entered "plain" (no explicit per-cell style carried over), matching
# the source workbook's last row. Column A/C still need the text-forcing
# trick (otherwise a date-shaped string silently becomes a date serial) -
# pulling blank/never-formatted cell format onto them afterwards keeps them
# on the sheet's implicit default style instead of picking up a new one.
$wsAmsin.Cells.Item(79, 1).NumberFormat = "@"
$wsAmsin.Cells.Item(79, 1).Value = "2023-04-12"
$wsAmsin.Cells.Item(1, 100).Copy()
$wsAmsin.Cells.Item(79, 1).PasteSpecial(-4122)

$wsAmsin.Cells.Item(78, 2).Copy()
$wsAmsin.Cells.Item(79, 2).PasteSpecial(-4122)
$wsAmsin.Cells.Item(79, 2).Value = 45028.61806041779

$wsAmsin.Cells.Item(79, 3).NumberFormat = "@"
$wsAmsin.Cells.Item(79, 3).Value = "176fstrtail"
$wsAmsin.Cells.Item(1, 100).Copy()
$wsAmsin.Cells.Item(79, 3).PasteSpecial(-4122)

$wsAmsin.Cells.Item(79, 4).Value = 89
$wsAmsin.Cells.Item(79, 5).Value = 87
$wsAmsin.Cells.Item(79, 6).Value = 2
$wsAmsin.Cells.Item(79, 7).Value = 2.78

# ---------------------------------------------------------------------
# BETA sheet: add rows 32-33
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-HistoryRow $wsBeta 32 "2023-03-13" 44998.54127167824 "174beta" 89 89 0 2.17 31
Add-HistoryRow $wsBeta 33 "2023-03-31" 45016.54783105324 "175beta" 89 89 0 2.17 32

# ---------------------------------------------------------------------
# AMS sheet: fix up row 41 (formatting + refreshed cached Run Time value),
# then add rows 42-45.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Range("A40:G40").Copy()
$wsAms.Range("A41:G41").PasteSpecial(-4122)
$wsAms.Cells.Item(41, 2).Value = 44977.8440737037

Add-HistoryRow $wsAms 42 "2023-03-01" 44986.69935061342 "173angularvrs" 89 89 0 2.34 41
Add-HistoryRow $wsAms 43 "2023-03-02" 44987.44013202546 "liveangular173" 89 89 0 2.33 42
Add-HistoryRow $wsAms 44 "2023-03-13" 44998.84198402778 "174live" 89 89 0 2.29 43
Add-HistoryRow $wsAms 45 "2023-03-31" 45016.81901835648 "175live" 89 89 0 2.21 44

Write-Host "Added AMSIN rows 76-79, BETA rows 32-33, AMS rows 42-45, refreshed AMS row 41."
